# Insert a new data row at sheet row 312 (pushing the existing rows
# 312-336 down to 313-337) and populate the new row with the new record.
# This mirrors the diff: dimension grows from A1:R336 to A1:R337 and a
# brand new "Apio" price record (fecha 2022-02-18) is inserted just after
# row 311.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(312).Insert()

$ws.Cells.Item(312, 1).Value = 3
$ws.Cells.Item(312, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(312, 3).Value = 'Coquimbo'
$ws.Cells.Item(312, 4).Value = 44610
$ws.Cells.Item(312, 5).Value = 5
$ws.Cells.Item(312, 6).Value = 100112017
$ws.Cells.Item(312, 7).Value = 'Apio'
$ws.Cells.Item(312, 8).Value = 'Americana (o)'
$ws.Cells.Item(312, 9).Value = 'Primera'
$ws.Cells.Item(312, 10).Value = 110
$ws.Cells.Item(312, 11).Value = 10000
$ws.Cells.Item(312, 12).Value = 10000
$ws.Cells.Item(312, 13).Value = 10000
$ws.Cells.Item(312, 14).Value = '$/docena de matas'
$ws.Cells.Item(312, 15).Value = 'Provincia de Santiago'
$ws.Cells.Item(312, 16).Value = 1667
$ws.Cells.Item(312, 17).Value = 6
$ws.Cells.Item(312, 18).Value = 'Hortaliza'
